$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.974.90'
$ws.Range('E2').Value = '  +0.32%  '
$ws.Range('D3').Value = '1.556.58'
$ws.Range('E3').Value = '  +0.74%  '
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '207.02'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.49%  '
$ws.Range('E6').Value = '  +0.36%  '
$ws.Range('E7').Value = '  -0.24%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.13'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.16%  '
$ws.Range('E9').Value = '  +0.19%  '
$ws.Range('E10').Value = '  +1.14%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0858'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.08%  '
$ws.Range('D12').Value = '1.778.10'
$ws.Range('E12').Value = '  +0.77%  '
$ws.Range('D13').Value = '1.555.88'
$ws.Range('E13').Value = '  +0.72%  '
$ws.Range('E14').Value = '  +1.45%  '
$ws.Range('E15').Value = '  +1.95%  '
$ws.Range('D16').Value = '26.983.75'
$ws.Range('E16').Value = '  +0.40%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.79'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.61%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '218.57'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.34%  '
$ws.Range('D19').Value = '0.0₃0697'
$ws.Range('E19').Value = '  +2.30%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.32'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.15%  '
$ws.Range('E21').Value = '  -0.26%  '
$ws.Range('E22').Value = '  +1.36%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.25'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.88%  '
$ws.Range('E24').Value = '  +0.76%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.50'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.36%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.66'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.00%  '
$ws.Range('E27').Value = '  +1.23%  '
$ws.Range('E28').Value = '  +1.24%  '
$ws.Range('E29').Value = '  -0.22%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0469'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E31').Value = '  -0.52%  '
$ws.Range('E32').Value = '  +0.83%  '
$ws.Range('D33').Value = '1.424.98'
$ws.Range('E33').Value = '  +5.10%  '
$ws.Range('E34').Value = '  +4.97%  '
$ws.Range('E35').Value = '  +3.82%  '
$ws.Range('E36').Value = '  +2.16%  '
$ws.Range('E37').Value = '  +0.20%  '
$ws.Range('E38').Value = '  +0.98%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.522'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.56%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.813'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.21%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.73'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.27%  '
$ws.Range('E42').Value = '  -0.21%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.33'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.987'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.25%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '64.59'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.10%  '
$ws.Range('E46').Value = '  +2.58%  '
$ws.Range('D47').Value = '1.691.76'
$ws.Range('E47').Value = '  +0.77%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '88.16'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.88%  '
$ws.Range('E49').Value = '  +1.94%  '
$ws.Range('D50').Value = '0.0₇0997'
$ws.Range('E50').Value = '  +2.15%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0958'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.31%  '
